$wb = $excel.ActiveWorkbook

$modelSheet = $wb.Worksheets.Item("Model")
$submodelsSheet = $wb.Worksheets.Item("Submodels")
$speciesTypesSheet = $wb.Worksheets.Item("Species types")
$rateLawsSheet = $wb.Worksheets.Item("Rate laws")
$parametersSheet = $wb.Worksheets.Item("Parameters")

# Simplify the rate law expressions on the Rate laws sheet
$rateLawsSheet.Range("C2").Value = "k1 * S1[c]"
$rateLawsSheet.Range("C3").Value = "k2 * S2[c] * S2[c]"

# Update the description text on the Model sheet
$modelSheet.Range("B2").Value = "Two reactions using two species in one compartment"

# Add a submodel column value to the Parameters sheet
$parametersSheet.Range("D3").Value = "ode_submodel"
$parametersSheet.Range("D4").Value = "ode_submodel"

# Update selections on a few sheets
$submodelsSheet.Range("A2").Select()
$speciesTypesSheet.Range("I2").Select()
$parametersSheet.Range("D3").Select()

# Make the Model sheet the active sheet/tab
$modelSheet.Activate()
